$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.905.20"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.821.32"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "707.37"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.29"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "3.819.55"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.38"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.42"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "4.466.23"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "3.823.28"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "70.917.85"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.31"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "494.38"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.61"
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.731"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.55"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.56"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "3.975.17"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.39"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.23"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.174"
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("D36").Value = "3.790.19"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.13"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.96"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.62"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "427.95"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.83"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.73"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.295"
$ws.Range("E51").Value = "  -2.40%  "
